# Apply crypto price/volume updates per commit "Updated cryptos list on Sat Aug 10 09:34:38 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text is a plain number would otherwise be auto-converted from text to a
# numeric value by Excel (losing formatting like trailing zeros, e.g. "10.60" -> 10.6).
# Force them to stay text: mark as Text format, set the value, then clear the transient
# number-format so the cell ends up with no explicit style (matching the source file).

$ws.Range("D2").Value = "60.747.43"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "2.616.53"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "513.98"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.25"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.40%  "
$ws.Range("E8").Value = "  -1.60%  "
$ws.Range("D9").Value = "2.625.75"
$ws.Range("E9").Value = "  -1.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.82"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +5.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.105"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("E13").Value = "  +1.76%  "
$ws.Range("D14").Value = "3.072.06"
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("D15").Value = "60.689.04"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.65"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000141"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("D18").Value = "2.619.46"
$ws.Range("E18").Value = "  -1.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.75"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "354.53"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.60"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.17"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.59%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.89"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.424"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.44%  "
$ws.Range("E26").Value = "  -0.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.995"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").Value = "0.0₃0844"
$ws.Range("E28").Value = "  -2.97%  "
$ws.Range("E29").Value = "  -3.80%  "
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.42"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.02"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.58"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.85"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.76%  "
$ws.Range("E35").Value = "  -2.12%  "
$ws.Range("E36").Value = "  -2.11%  "
$ws.Range("E37").Value = "  +3.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.49"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.36"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.844"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.77%  "
$ws.Range("E41").Value = "  -0.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "293.64"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -6.13%  "
$ws.Range("E43").Value = "  +0.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.624"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.995"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0554"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.81"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.93"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("E49").Value = "  -1.27%  "
$ws.Range("E50").Value = "  +0.40%  "
$ws.Range("D51").Value = "2.002.09"
$ws.Range("E51").Value = "  -2.72%  "
